$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.915.99"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.901.97"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'0.7668"
$ws.Range("E5").Value = "  +5.37%  "
$ws.Range("D6").Value = "'240.58"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "'0.3069"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'25.65"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").Value = "'0.06854"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "'0.07976"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.915.86"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "'0.7452"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "'5.170"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "'91.17"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "29.905.64"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'13.98"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'5.964"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "'243.34"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'0.000007704"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'6.959"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").Value = "'166.77"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'9.251"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'18.74"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("D28").Value = "'2.048"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").Value = "'1.411"
$ws.Range("E29").Value = "  +4.78%  "
$ws.Range("D30").Value = "'1.517"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "'4.261"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "'4.089"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'0.05263"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").Value = "'0.7290"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'2.713"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'0.01930"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "'2.773"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'6.181"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").Value = "'0.4422"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'72.19"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "'1.890"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("D44").Value = "'0.8295"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").Value = "'7.632"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "'100.16"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'9.776"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "2.052.04"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("D49").Value = "'36.16"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.477"
$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05941"
$ws.Range("E51").Value = "  +0.10%  "
